$wb = $excel.ActiveWorkbook

# --- HH sheet updates ---
$hh = $wb.Worksheets.Item("HH")

# "valor HH" rate increased -> update the "valor HH total" formula in S3
$hh.Range("S3").Formula = "=(284.03)*S2"

# Newly logged work hours for R9:R11 (same number format/style as R8, [h]:mm:ss elapsed time)
$hh.Range("R9").NumberFormat = "[h]:mm:ss"
$hh.Range("R9").Value = 0.98472222222222217

$hh.Range("R10").NumberFormat = "[h]:mm:ss"
$hh.Range("R10").Value = 0.84305555555555556

$hh.Range("R11").NumberFormat = "[h]:mm:ss"
$hh.Range("R11").Value = 0.93958333333333333

# R20 carries the time-format (h:mm) used further down the same column
$hh.Range("R20").NumberFormat = "h:mm"

# Running total formula now includes the newly added hours
$hh.Range("Q6").Formula = "= Q8+R8+R9+R10+R11"

# Recalculate so all dependent formulas (Presupuesto sheet included) refresh
$excel.Calculate()

# Last active selection left on the HH sheet
$hh.Range("S3").Select()

# --- Presupuesto sheet updates ---
$pres = $wb.Worksheets.Item("Presupuesto")
$pres.Activate()
$pres.Range("B11").Select()
